$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "paymentReference" column (E) - header + 13 data rows.
# Must be written fully before column F so new shared strings are
# appended in the same order as the target workbook (paymentReference*
# entries precede dealerName* entries in sharedStrings.xml).
$ws.Range("E1").Value = "paymentReference"
$ws.Range("E2").Value = "paymentReference1"
$ws.Range("E3").Value = "paymentReference2"
$ws.Range("E4").Value = "paymentReference3"
$ws.Range("E5").Value = "paymentReference4"
$ws.Range("E6").Value = "paymentReference5"
$ws.Range("E7").Value = "paymentReference6"
$ws.Range("E8").Value = "paymentReference7"
$ws.Range("E9").Value = "paymentReference8"
$ws.Range("E10").Value = "paymentReference9"
$ws.Range("E11").Value = "paymentReference10"
$ws.Range("E12").Value = "paymentReference11"
$ws.Range("E13").Value = "paymentReference12"
$ws.Range("E14").Value = "paymentReference13"

# Add "dealerName" column (F) - header + 13 data rows.
$ws.Range("F1").Value = "dealerName"
$ws.Range("F2").Value = "dealerName1"
$ws.Range("F3").Value = "dealerName2"
$ws.Range("F4").Value = "dealerName3"
$ws.Range("F5").Value = "dealerName4"
$ws.Range("F6").Value = "dealerName5"
$ws.Range("F7").Value = "dealerName6"
$ws.Range("F8").Value = "dealerName7"
$ws.Range("F9").Value = "dealerName8"
$ws.Range("F10").Value = "dealerName9"
$ws.Range("F11").Value = "dealerName10"
$ws.Range("F12").Value = "dealerName11"
$ws.Range("F13").Value = "dealerName12"
$ws.Range("F14").Value = "dealerName13"

# Update the selected cell to match the target workbook's active cell.
$ws.Range("F8").Select()
